$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width
$ws.Columns.Item(1).ColumnWidth = 16.42578125

# Update cell values in column A
$values = @(
    0.12926056049198564,
    -0.0099999991685244538,
    -0.067712871276739861,
    0.28398915160247284,
    -0.0059999991904930994,
    -0.0059999991577299738,
    -0.019999999016093284,
    -0.019999999004613578,
    -0.0059999991284431786,
    -0.0059999991147563492,
    -0.0044999991276739593,
    -0.0059999991097812178,
    -0.0059999990960539762,
    -0.011999999033037057,
    -0.0059999990889076926,
    -0.005999999085952723,
    -0.0059999990821086868,
    -0.0089999990517437567,
    -0.0089999991780547184,
    -0.0089999991711486871,
    -0.008999999170162809,
    -0.0089999991694398318,
    -0.0089999991639420074,
    -0.041999998826843132,
    -0.041999998820732465,
    -0.0059999991550796494,
    -0.00090999155841853252,
    -0.0059999991365744521,
    -0.011999999067482392,
    -0.019999998984197465,
    -0.0068026635077593767,
    -0.020999998966333422,
    -0.0059999991135226693
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
